$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 227-229 continue the daily data series. Copy the formatting
# (date style, borders, number format) from the last existing row (226)
# down into the three new rows before writing their values.
$ws.Range("A226").Copy()
$ws.Range("A227:A229").PasteSpecial(-4122)

$ws.Range("A227").Value = 44301
$ws.Range("B227").Value = 1
$ws.Range("C227").Value = 4
$ws.Range("D227").Value = 175.1313485113835

$ws.Range("A228").Value = 44302
$ws.Range("B228").Value = 0
$ws.Range("C228").Value = 4
$ws.Range("D228").Value = 175.1313485113835

$ws.Range("A229").Value = 44303
$ws.Range("B229").Value = 0
$ws.Range("C229").Value = 4
$ws.Range("D229").Value = 175.1313485113835
